$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 80
$ws.Cells.Item(2, 3).Value = "dog/dog025.jpg"
$ws.Cells.Item(2, 4).Value = "proben"
$ws.Cells.Item(2, 5).Value = "dog"
$ws.Cells.Item(3, 2).Value = 92
$ws.Cells.Item(3, 3).Value = "house/house025.jpg"
$ws.Cells.Item(3, 4).Value = "stoppen"
$ws.Cells.Item(3, 5).Value = "house"
$ws.Cells.Item(4, 2).Value = 121
$ws.Cells.Item(4, 3).Value = "house/house006.jpg"
$ws.Cells.Item(4, 4).Value = "wählen"
$ws.Cells.Item(4, 5).Value = "house"
$ws.Cells.Item(5, 2).Value = 114
$ws.Cells.Item(5, 3).Value = "dog/dog024.jpg"
$ws.Cells.Item(5, 4).Value = "kranken"
$ws.Cells.Item(5, 5).Value = "dog"
$ws.Cells.Item(6, 2).Value = 7
$ws.Cells.Item(6, 3).Value = "house/house019.jpg"
$ws.Cells.Item(6, 4).Value = "hacken"
$ws.Cells.Item(6, 5).Value = "house"
$ws.Cells.Item(7, 2).Value = 13
$ws.Cells.Item(7, 3).Value = "house/house024.jpg"
$ws.Cells.Item(7, 4).Value = "münzen"
$ws.Cells.Item(7, 5).Value = "house"
$ws.Cells.Item(8, 2).Value = 122
$ws.Cells.Item(8, 3).Value = "house/house000.jpg"
$ws.Cells.Item(8, 4).Value = "bauen"
$ws.Cells.Item(8, 5).Value = "house"
$ws.Cells.Item(9, 2).Value = 81
$ws.Cells.Item(9, 3).Value = "dog/dog013.jpg"
$ws.Cells.Item(9, 4).Value = "fließen"
$ws.Cells.Item(9, 5).Value = "dog"
$ws.Cells.Item(10, 2).Value = 42
$ws.Cells.Item(10, 3).Value = "house/house020.jpg"
$ws.Cells.Item(10, 4).Value = "lügen"
$ws.Cells.Item(10, 5).Value = "house"
$ws.Cells.Item(11, 2).Value = 99
$ws.Cells.Item(11, 3).Value = "dog/dog003.jpg"
$ws.Cells.Item(11, 4).Value = "lächeln"
$ws.Cells.Item(11, 5).Value = "dog"
$ws.Cells.Item(12, 2).Value = 115
$ws.Cells.Item(12, 3).Value = "dog/dog028.jpg"
$ws.Cells.Item(12, 4).Value = "legen"
$ws.Cells.Item(12, 5).Value = "dog"
$ws.Cells.Item(13, 2).Value = 95
$ws.Cells.Item(13, 3).Value = "dog/dog021.jpg"
$ws.Cells.Item(13, 4).Value = "mögen"
$ws.Cells.Item(13, 5).Value = "dog"
$ws.Cells.Item(14, 2).Value = 78
$ws.Cells.Item(14, 3).Value = "dog/dog030.jpg"
$ws.Cells.Item(14, 4).Value = "leeren"
$ws.Cells.Item(14, 5).Value = "dog"
$ws.Cells.Item(15, 2).Value = 46
$ws.Cells.Item(15, 3).Value = "house/house001.jpg"
$ws.Cells.Item(15, 4).Value = "planen"
$ws.Cells.Item(15, 5).Value = "house"
$ws.Cells.Item(16, 2).Value = 21
$ws.Cells.Item(16, 3).Value = "house/house010.jpg"
$ws.Cells.Item(16, 4).Value = "narren"
$ws.Cells.Item(16, 5).Value = "house"
$ws.Cells.Item(17, 2).Value = 61
$ws.Cells.Item(17, 3).Value = "dog/dog023.jpg"
$ws.Cells.Item(17, 4).Value = "wachsen"
$ws.Cells.Item(17, 5).Value = "dog"
$ws.Cells.Item(18, 2).Value = 93
$ws.Cells.Item(18, 3).Value = "house/house005.jpg"
$ws.Cells.Item(18, 4).Value = "duschen"
$ws.Cells.Item(18, 5).Value = "house"
$ws.Cells.Item(19, 2).Value = 32
$ws.Cells.Item(19, 3).Value = "dog/dog005.jpg"
$ws.Cells.Item(19, 4).Value = "ärgern"
$ws.Cells.Item(19, 5).Value = "dog"
$ws.Cells.Item(20, 2).Value = 25
$ws.Cells.Item(20, 3).Value = "dog/dog002.jpg"
$ws.Cells.Item(20, 4).Value = "piepen"
$ws.Cells.Item(20, 5).Value = "dog"
$ws.Cells.Item(21, 2).Value = 62
$ws.Cells.Item(21, 3).Value = "dog/dog008.jpg"
$ws.Cells.Item(21, 4).Value = "betteln"
$ws.Cells.Item(21, 5).Value = "dog"
$ws.Cells.Item(22, 2).Value = 44
$ws.Cells.Item(22, 3).Value = "dog/dog006.jpg"
$ws.Cells.Item(22, 4).Value = "nullen"
$ws.Cells.Item(22, 5).Value = "dog"
$ws.Cells.Item(23, 2).Value = 94
$ws.Cells.Item(23, 3).Value = "house/house014.jpg"
$ws.Cells.Item(23, 4).Value = "küssen"
$ws.Cells.Item(23, 5).Value = "house"
$ws.Cells.Item(24, 2).Value = 19
$ws.Cells.Item(24, 3).Value = "dog/dog022.jpg"
$ws.Cells.Item(24, 4).Value = "wehen"
$ws.Cells.Item(24, 5).Value = "dog"
$ws.Cells.Item(25, 2).Value = 117
$ws.Cells.Item(25, 3).Value = "house/house003.jpg"
$ws.Cells.Item(25, 4).Value = "süßen"
$ws.Cells.Item(25, 5).Value = "house"
$ws.Cells.Item(26, 2).Value = 60
$ws.Cells.Item(26, 3).Value = "house/house008.jpg"
$ws.Cells.Item(26, 4).Value = "öffnen"
$ws.Cells.Item(26, 5).Value = "house"
$ws.Cells.Item(27, 2).Value = 36
$ws.Cells.Item(27, 3).Value = "house/house012.jpg"
$ws.Cells.Item(27, 4).Value = "spenden"
$ws.Cells.Item(27, 5).Value = "house"
$ws.Cells.Item(28, 2).Value = 28
$ws.Cells.Item(28, 3).Value = "house/house007.jpg"
$ws.Cells.Item(28, 4).Value = "trotzen"
$ws.Cells.Item(28, 5).Value = "house"
$ws.Cells.Item(29, 2).Value = 63
$ws.Cells.Item(29, 3).Value = "house/house002.jpg"
$ws.Cells.Item(29, 4).Value = "dienen"
$ws.Cells.Item(29, 5).Value = "house"
$ws.Cells.Item(30, 2).Value = 105
$ws.Cells.Item(30, 3).Value = "house/house009.jpg"
$ws.Cells.Item(30, 4).Value = "holen"
$ws.Cells.Item(30, 5).Value = "house"
$ws.Cells.Item(31, 2).Value = 75
$ws.Cells.Item(31, 3).Value = "dog/dog016.jpg"
$ws.Cells.Item(31, 4).Value = "hassen"
$ws.Cells.Item(31, 5).Value = "dog"
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = "dog/dog020.jpg"
$ws.Cells.Item(32, 4).Value = "frischen"
$ws.Cells.Item(32, 5).Value = "dog"
$ws.Cells.Item(33, 2).Value = 126
$ws.Cells.Item(33, 3).Value = "dog/dog001.jpg"
$ws.Cells.Item(33, 4).Value = "achten"
$ws.Cells.Item(33, 5).Value = "dog"
